$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.988.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.898.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.8334'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.82'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3287'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.58'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07048'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08079'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7612'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.895.88'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.250'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.987.32'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.10'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.866'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007761'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.154.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.966'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1726'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +25.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.259'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.095'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.363'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.514'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05956'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.288'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.078'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.269'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7313'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.725'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01916'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.779'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4442'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.51'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.857'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.15%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8537'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.06%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.892'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.00'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.777'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '994.59'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.046.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.04%  '
